$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.077.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "'3.563.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'606.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'145.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "

$ws.Range("D7").Value = "'3.562.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.48%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  +3.58%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").Value = "'7.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.02%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").Value = "'4.166.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").Value = "'29.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'3.568.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Value = "'66.185.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.82%  "

$ws.Range("D20").Value = "'6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").Value = "'14.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").Value = "'429.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("E23").Value = "  +4.65%  "

$ws.Range("D24").Value = "'79.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("D25").Value = "'3.705.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.56%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("D29").Value = "'7.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").Value = "'9.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "'25.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.93%  "

$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").Value = "'3.556.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.62%  "

$ws.Range("E35").Value = "  -6.02%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("D38").Value = "'7.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.15%  "

$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "'175.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.29%  "

$ws.Range("D42").Value = "'0.0847"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "

$ws.Range("E43").Value = "  +2.38%  "

$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("D46").Value = "'46.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "

$ws.Range("D47").Value = "'25.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").Value = "'2.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "

$ws.Range("D50").Value = "'23.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.28%  "

$ws.Range("E51").Value = "  +0.24%  "
